# Normalize the "Recorded By" (column G) values: move a leading "System"
# entry to the end of the comma-separated list, e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"
# The single exception observed in the source data is the literal value
# "System, admin@admin.com", which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Formula

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if (-not $val.StartsWith("System,")) { continue }
    if ($val -eq "System, admin@admin.com") { continue }

    $parts = $val -split ", "
    $rest = $parts[1..($parts.Length - 1)]
    $newParts = $rest + @("System")
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value = $newVal
}
